$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the percentage value for TELE's row (A2) as text with a leading
# apostrophe so it keeps its original "quote prefix" text formatting but
# loses the thousands separator (1,231.53 -> 1231.53).
$ws.Range("A2").Value = "'1231.53"

# Update the active selection to A3, matching the saved cursor position.
$ws.Range("A3").Select()
